# "updated order estimate, fixed minor schematic opcode oopsies"
#
# The "notes" worksheet (a small running tally of parts still needed per
# sub-board) gets reworked: a new column D ("x" = still need to order /
# "no" = already covered) is added throughout, a couple of part
# descriptions are corrected (resistor footprints etc.), a new "pico
# header" line is inserted into the memory-unit block, the "other:" and
# "mainboard:" sections gain several new line items, and the whole table
# grows from A1:C28 to A1:D34.
#
# Rewriting the table from scratch (clear, then re-enter every cell) is
# simpler and less error-prone than trying to replay the individual
# row-insert operations the author performed interactively.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("notes")

$ws.Cells.Clear()

function Set-Cells($row, $a, $b, $c, $d) {
    if ($null -ne $a) { $ws.Cells.Item($row, 1).Value = $a }
    if ($null -ne $b) { $ws.Cells.Item($row, 2).Value = $b }
    if ($null -ne $c) { $ws.Cells.Item($row, 3).Value = $c }
    if ($null -ne $d) { $ws.Cells.Item($row, 4).Value = $d }
}

# mem units (2):
Set-Cells  1 "mem units (2):" "m header 2x21"     "4 (5)"        "x"
Set-Cells  2 $null             "IC sockets pic"    "20 (30)"      "x"
Set-Cells  3 $null             "pico `"header`""   "2 (2)"        "x"
Set-Cells  4 $null             "variety m header"  "4pin (50pin)" "x"
Set-Cells  5 $null             "resistors 0805!!!" "book"         "no"
Set-Cells  6 $null             "capacitors 0805"   "-"            "no"
Set-Cells  7 $null             "schottky"          "(?)"          "no"

# IO units (4):
Set-Cells  9 "IO units (4):"  "m header 2x7"      "4 (5)"        "x"
Set-Cells 10 $null             "IC sockets pic"    "4 (30)"       "x"
Set-Cells 11 $null             "variety m header"  "8pin (50pin)" "x"
Set-Cells 12 $null             "resistors 0805!!!" "book"         "no"
Set-Cells 13 $null             "capacitors 0805"   "-"            "no"

# sequencer (2):
Set-Cells 15 "sequencer (2):" "m header 2x5"      "4 (5)"        "x"
Set-Cells 16 $null             "IC sockets pic"    "2 (30)"       "x"
Set-Cells 17 $null             "variety m header"  "4pin (50pin)" "x"
Set-Cells 18 $null             "resistors 0805!!!" "book"         "no"
Set-Cells 19 $null             "capacitors 0805"   "-"            "no"
Set-Cells 20 $null             "buttons"           10             "x"
Set-Cells 21 $null             "switches"          10             "x"

# other:
Set-Cells 23 "other:"         "breadboards"       "10?"          "unavailable"
Set-Cells 24 $null             "IDC m header 2x4"  4              "x"
Set-Cells 25 $null             "on/off switch"     2              "x"
Set-Cells 26 $null             "power switch"      2              "x"

# mainboard:
Set-Cells 28 "mainboard:"     "resistors 0603 10k" "100?"        "x"
Set-Cells 29 $null             "capacitors 0805"   "-"            "no"
Set-Cells 30 $null             "full adder"        20             "x"
Set-Cells 31 $null             "multiplexer"       5              "x"
Set-Cells 32 $null             "f header 2x21"     "4 (5)"        "x"
Set-Cells 33 $null             "f header 2x5"      "2 (5)"        "x"
Set-Cells 34 $null             "f header 2x7"      "4 (5)"        "x"
